# Applies the "Structure changes, improve documentation." commit:
#  1. HEADER sheet: relabel rows 3-6 (DOMAIN/CATEGORY/SOURCE_ORG/SOURCE_PERSON),
#     clearing the values that used to belong to CATEGORY/SUB_CATEGORY.
#  2. SHARE_ sheet: insert a new "-" entry at the top of column A (rows shift
#     down by one; column B is untouched).
#  3. DIV_ sheet: insert a new "-" entry at the top of columns A and B
#     (rows shift down by one; columns C/D are untouched).
#  4. SHARE sheet: widen the SHARE_ column-A list-validation range to match
#     the extra row (A1:A34 -> A1:A35).
#  5. DIV sheet: widen the DIV_ column-B list-validation range to match the
#     extra row (B1:B14 -> B1:B15).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. HEADER
# ---------------------------------------------------------------------
$header = $wb.Worksheets.Item("HEADER")

$header.Range("A3").Value = "DOMAIN"
$header.Range("A4").Value = "CATEGORY"
$header.Range("A5").Value = "SOURCE_ORG"
$header.Range("B5").Value = ""
$header.Range("A6").Value = "SOURCE_PERSON"
$header.Range("B6").Value = ""

# ---------------------------------------------------------------------
# 2. SHARE_ : prepend "-" to column A, shifting the rest down one row
# ---------------------------------------------------------------------
$shareHidden = $wb.Worksheets.Item("SHARE_")

$shareColA = @("A","B","D","H","M","N","Q","S","W","A2","A3","A4","A5","A10","A20","A30","A_3","M2","M_2","M_3","W2","W3","W4","W_2","W_3","D_2","H2","H3","I","OA","OM","_O","_U","_Z")

$shareColA = @("-") + $shareColA

for ($i = 0; $i -lt $shareColA.Length; $i++) {
    $shareHidden.Cells.Item($i + 1, 1).Value = $shareColA[$i]
}

# ---------------------------------------------------------------------
# 3. DIV_ : prepend "-" to columns A and B, shifting the rest down one row
# ---------------------------------------------------------------------
$divHidden = $wb.Worksheets.Item("DIV_")

$divColA = @("AN","SA","QA","BM","MO","WE","DA")
$divColA = @("-") + $divColA

$divColB = @("10","20","11","12","13","21","22","23","31","32","33","91","92","93")
$divColB = @("-") + $divColB

for ($i = 0; $i -lt $divColA.Length; $i++) {
    $divHidden.Cells.Item($i + 1, 1).Value = $divColA[$i]
}

for ($i = 0; $i -lt $divColB.Length; $i++) {
    $divHidden.Cells.Item($i + 1, 2).Value = $divColB[$i]
}

# ---------------------------------------------------------------------
# 4. SHARE : widen the list validation that points at SHARE_ column A
# ---------------------------------------------------------------------
$share = $wb.Worksheets.Item("SHARE")
$share.Range("C4:C20").Validation.Formula1 = "'SHARE_'!`$A`$1:`$A`$35"

# ---------------------------------------------------------------------
# 5. DIV : widen the list validation that points at DIV_ column B
# ---------------------------------------------------------------------
$div = $wb.Worksheets.Item("DIV")
$div.Range("G4:G20").Validation.Formula1 = "'DIV_'!`$B`$1:`$B`$15"
